# Merge - Opp Test Data, ENg Detail, Add Counterparty - 10 Oct 2025
#
# 1) On the "Users" sheet, the CaoUser test record "Liz Hedgcock" is
#    replaced with "Blaise Brunda".
# 2) The "Users" sheet is reordered to be the first (left-most) tab in
#    the workbook (it used to sit after AddOpportunity/AppName/ModuleName).

$wb = $excel.ActiveWorkbook

# --- 1. Update the counterparty/user name on the Users sheet -----------
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("B2").Value = "Blaise Brunda"

# --- 2. Make Users the active sheet with its own remembered selection --
#     (do this before moving, since Move() can invalidate the handle)
$usersSheet.Activate()
$usersSheet.Range("F7").Select()

# --- 3. Move Users to be the very first worksheet tab -------------------
$usersSheet.Move($wb.Worksheets.Item(1))

# --- 4. Re-activate it by a fresh lookup so it stays the shown tab ------
$frontUsersSheet = $wb.Worksheets.Item("Users")
$frontUsersSheet.Activate()
